# Add the "FormulaeTypes" worksheet right after the existing "ValueTypes"
# sheet, demonstrating simple / shared / array formulae.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)
$ws2.Name = "FormulaeTypes"

# --- Column A labels -------------------------------------------------
# Written in this particular order so that the workbook's shared-string
# table ends up holding the strings in the same order as the reference
# workbook: Hello, Simple, Shared, Sharing, Array (single),
# Arraying (multiple), Arrayed (multiple).
$ws2.Range("A1").Value = "Simple"
$ws2.Range("A3").Value = "Shared"
$ws2.Range("A4").Value = "Shared"
$ws2.Range("A2").Value = "Sharing"
$ws2.Range("A5").Value = "Array (single)"
$ws2.Range("A6").Value = "Arraying (multiple)"
$ws2.Range("A7").Value = "Arrayed (multiple)"
$ws2.Range("A8").Value = "Arrayed (multiple)"

# --- Column B formulae ------------------------------------------------
# Simple formula.
$ws2.Range("B1").Formula = "=1+1"

# A formula that is identical to, but not shared with, the one below.
$ws2.Range("B2").Formula = "=COSH(2*PI())"

# A shared formula spanning B3:B4 (assign to the whole range at once so
# the engine emits a single master <f t="shared"> cell plus a follower).
$ws2.Range("B3:B4").Formula = "=COSH(2*PI())"

# A single-cell array formula.
$ws2.Range("B5").FormulaArray = "=B1:B4"

# A multi-cell array formula spilling across B6:B8.
$ws2.Range("B6:B8").FormulaArray = '=IF(B3:B5=8,"Eight","Not Eight")'

# --- Cosmetics to mirror the reference workbook -----------------------
# Column A sized (≈16.5 "chars" once Excel's standard padding is added).
$ws2.Columns.Item(1).ColumnWidth = 15.6666666667

# Page margins (inches, as shown in Page Setup); COM takes points.
$ws2.PageSetup.LeftMargin   = 0.75 * 72
$ws2.PageSetup.RightMargin  = 0.75 * 72
$ws2.PageSetup.TopMargin    = 1    * 72
$ws2.PageSetup.BottomMargin = 1    * 72
$ws2.PageSetup.HeaderMargin = 0.5  * 72
$ws2.PageSetup.FooterMargin = 0.5  * 72
$ws2.PageSetup.PaperSize    = 9     # matches the reference sheet
$ws2.PageSetup.Orientation  = 1     # xlPortrait

# Selection / active sheet: FormulaeTypes becomes the active tab, with
# B5 selected - matching the reference workbook's saved UI state.
$ws2.Range("B5").Select()
$ws2.Activate()
